# add save column in s_vals sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Save" column, matching the style of the other headers (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values for rows 2-8
$saveValues = @(0, 0, 1, 1, 1, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
